$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) "There are 3 levels ..." paragraph: a trailing "." is typed as a new run
#    at the end of the sentence.
# ---------------------------------------------------------------------------
$introPara = $d.Paragraphs(2)
$introRange = $introPara.Range
$introXml = "<w:p $wNs w14:paraId=`"24C5E203`" w14:textId=`"34D7B8C1`" w:rsidR=`"00805F81`" w:rsidRDefault=`"00805F81`" w:rsidP=`"00805F81`">" +
            "<w:r><w:t>There are 3 levels of results, Green (lowest), Yellow, and Red (highest)</w:t></w:r>" +
            "<w:r><w:t>.</w:t></w:r>" +
            "</w:p>"
$introRange.InsertXML($introXml)

# ---------------------------------------------------------------------------
# 2) "2 different authors once the students commits begin" paragraph gets
#    proofing marks around "commits" (subject/verb agreement) and is
#    followed by a brand new bullet (one level deeper) that continues the
#    thought about ACES looking for the author "Default". The _GoBack
#    bookmark that used to sit at the end of the first paragraph now sits
#    inside the new paragraph.
# ---------------------------------------------------------------------------
$authorsPara = $d.Paragraphs(8)
$authorsRange = $authorsPara.Range
$authorsXml = "<w:p $wNs w14:paraId=`"262B329A`" w14:textId=`"31B210FE`" w:rsidR=`"00805F81`" w:rsidRDefault=`"00805F81`" w:rsidP=`"00527BC6`">" +
              "<w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr>" +
              "<w:r><w:t xml:space=`"preserve`">2 different authors once the students </w:t></w:r>" +
              "<w:proofErr w:type=`"gramStart`"/>" +
              "<w:r><w:t>commits</w:t></w:r>" +
              "<w:proofErr w:type=`"gramEnd`"/>" +
              "<w:r><w:t xml:space=`"preserve`"> begi</w:t></w:r>" +
              "<w:r><w:t>n</w:t></w:r>" +
              "</w:p>" +
              "<w:p $wNs>" +
              "<w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"1`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr>" +
              "<w:r><w:t xml:space=`"preserve`">Right </w:t></w:r>" +
              "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/>" +
              "<w:bookmarkEnd w:id=`"0`"/>" +
              "<w:r><w:t xml:space=`"preserve`">now, ACES </w:t></w:r>" +
              "<w:proofErr w:type=`"gramStart`"/>" +
              "<w:r><w:t>looks</w:t></w:r>" +
              "<w:proofErr w:type=`"gramEnd`"/>" +
              "<w:r><w:t xml:space=`"preserve`"> for the author &#8220;Default&#8221;. This would be good to change, but will require a fair bit of refactoring</w:t></w:r>" +
              "</w:p>"
$authorsRange.InsertXML($authorsXml)
